# Anonymize "fedcore" -> "approach" and add top/bottom (and top/bottom/right)
# thin borders to the merged header cells C1:D1 (and F1:G1 on sheet 2),
# matching the "change" header cell style already used elsewhere.
# Also clears a stray empty inline-string cell (G5) on sheet 2.

$wb = $excel.ActiveWorkbook

# xlLineStyle / xlBorderWeight constants used below
$xlNone = -4142
$xlContinuous = 1
$xlThin = 2

function Set-HeaderBorderStyle {
    param(
        [__ComObject]$ws,
        [string]$cellAddr,
        [bool]$rightThin
    )

    $dst = $ws.Range($cellAddr)

    # Start from a cell that already carries the default (unstyled) format
    # so the new style reuses font 0 / fill 0 instead of the bold header font.
    $src = $ws.Range("B4")
    $src.Copy()
    $dst.PasteSpecial(-4122) # xlPasteFormats

    $dst.Borders.Item(7).LineStyle = $xlNone        # xlEdgeLeft
    $dst.Borders.Item(8).LineStyle = $xlContinuous   # xlEdgeTop
    $dst.Borders.Item(9).LineStyle = $xlContinuous   # xlEdgeBottom
    if ($rightThin) {
        $dst.Borders.Item(10).LineStyle = $xlContinuous # xlEdgeRight
    } else {
        $dst.Borders.Item(10).LineStyle = $xlNone
    }
}

# --- Sheet 1: quality_comparison -------------------------------------------
$ws1 = $wb.Worksheets.Item("quality_comparison")

Set-HeaderBorderStyle -ws $ws1 -cellAddr "C1" -rightThin $false
Set-HeaderBorderStyle -ws $ws1 -cellAddr "D1" -rightThin $true

$ws1.Range("C2").Value = "approach"

# --- Sheet 2: computational_comparison --------------------------------------
$ws2 = $wb.Worksheets.Item("computational_comparison")

Set-HeaderBorderStyle -ws $ws2 -cellAddr "C1" -rightThin $false
Set-HeaderBorderStyle -ws $ws2 -cellAddr "D1" -rightThin $true
Set-HeaderBorderStyle -ws $ws2 -cellAddr "F1" -rightThin $false
Set-HeaderBorderStyle -ws $ws2 -cellAddr "G1" -rightThin $true

$ws2.Range("C2").Value = "approach"
$ws2.Range("F2").Value = "approach"

# Drop the stray empty inline-string cell entirely (matches diff removal).
$ws2.Range("G5").ClearContents()
